# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") on rows 2-30, recomputed after switching the
# source statistic from Strike# to K.
$newK = @{
    2  = 3
    3  = 5
    4  = 3
    5  = 4
    6  = 2
    7  = 1
    8  = 4
    9  = 3
    10 = 3
    11 = 2
    12 = 7
    13 = 4
    14 = 4
    15 = 5
    16 = 5
    17 = 6
    18 = 3
    19 = 3
    20 = 6
    21 = 3
    22 = 10
    23 = 5
    24 = 6
    25 = 2
    26 = 5
    27 = 7
    28 = 8
    29 = 2
    30 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
